# WW101-Binder-Cover.pptx edit
#  1) Bump the cached "datetimeFigureOut" date field (slide master + all
#     slide layouts) from 7/3/2018 -> 8/7/2018.
#  2) Slide 1 "TextBox 6" (part/rev stamp): reposition/resize slightly and
#     bump the revision text 002-23599 *A -> 002-23599 *B.

$p = $ppt.ActivePresentation

# --- 1) Date placeholders --------------------------------------------------
$oldDate = "7/3/2018"
$newDate = "8/7/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# --- 2) Slide 1 revision stamp ---------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "002-23599 *A") {
        $shp.Left = 455.9544881889764
        $shp.Top = 692.4524409448819
        $shp.Width = 74.74795275590552

        $tr = $shp.TextFrame.TextRange
        $tr.Delete()
        $null = $tr.InsertAfter("002-23599 *B")
    }
}
